$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("D9")
$r.Font.Underline = 2
$r.Font.Color = 16711680
Write-Host "COLOR: " $r.Font.Color
Write-Host "UNDERLINE: " $r.Font.Underline
